$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.556.11'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '1.752.34'
$ws.Range("E3").Value = '  -3.42%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4466'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3613'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07487'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.106'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.98%  '
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.033'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.177'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.18%  '
$ws.Range("D16").Value = '1.752.32'
$ws.Range("E16").Value = '  -3.85%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.66%  '
$ws.Range("E18").Value = '  -1.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06416'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.03'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.842'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.53%  '
$ws.Range("D23").Value = '27.598.43'
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '1.956.19'
$ws.Range("E28").Value = '  -3.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.128'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.088'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.52%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09024'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.639'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.540'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02302'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2100'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6369'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05957'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.04%  '
$ws.Range("E40").Value = '  -4.89%  '
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.389'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.782'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.717'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5878'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.960'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.155'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06852'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.05%  '
